$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that sits right after the
#    "Programmer" run in the Work Experience section. Word will
#    renumber the remaining bookmark ids automatically (5->4, 6->5,
#    7->6, 8 stays 8 for now until step 3 reinserts a _GoBack later
#    in the doc).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Education section: turn "B.Sc. in Computer Science" into
#    "Bachelor of Computing" and make sure the surrounding text keeps
#    the same run boundaries as before ("201" / "3" / " " / the
#    replaced phrase / ", University of " / "Brasília").
#    We drop temporary bookmarks right at the boundaries first (that
#    does not disturb the existing runs) so that the later text
#    replace only touches the isolated run instead of merging with
#    its neighbours.
# ------------------------------------------------------------------
$match = $d.Content
$match.Find.Execute("B.Sc. in Computer Science")
$matchStart = $match.Start
$matchEnd = $match.End

$splitBefore = $d.Range($matchStart, $matchStart)
$d.Bookmarks.Add("_ZZtmpSplitBefore", $splitBefore)

$splitAfter = $d.Range($matchEnd, $matchEnd)
$d.Bookmarks.Add("_ZZtmpSplitAfter", $splitAfter)

$replaceRange = $d.Content
$replaceRange.Find.Execute("B.Sc. in Computer Science")
$replaceRange.Text = "Bachelor of Computing"

# Remove the leading helper bookmark - it leaves the run split as-is.
$d.Bookmarks("_ZZtmpSplitBefore").Delete()

# ------------------------------------------------------------------
# 3. Turn the trailing helper bookmark into the new "_GoBack" bookmark
#    (matches Word re-recording the last edit position right after
#    typing "Bachelor of Computing").
# ------------------------------------------------------------------
$afterBm = $d.Bookmarks("_ZZtmpSplitAfter")
$gbRange = $d.Range($afterBm.Start, $afterBm.End)
$afterBm.Delete()
$d.Bookmarks.Add("_GoBack", $gbRange)

# ------------------------------------------------------------------
# 4. "Default Paragraph Font" character style should no longer be
#    semi-hidden (it became visible once it was actually put to use).
# ------------------------------------------------------------------
$dpf = $d.Styles("Default Paragraph Font")
$dpf.UnhideWhenUsed = $true
